# Append new daily COVID cumulative-death rows (2021-02-15 .. 2021-03-07)
# to Sheet1, then leave the selection/view positioned on the new last row,
# matching the state after the data update was pasted in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each row: Date(serial), DeathCovid, DeathWithCovid, Total
$newRows = @(
    @(44242, 6063, 1266, 7329),
    @(44243, 6168, 1273, 7441),
    @(44244, 6271, 1288, 7559),
    @(44245, 6350, 1292, 7642),
    @(44246, 6424, 1301, 7725),
    @(44247, 6505, 1313, 7818),
    @(44248, 6577, 1334, 7911),
    @(44249, 6671, 1343, 8014),
    @(44250, 6775, 1360, 8135),
    @(44251, 6859, 1370, 8229),
    @(44252, 6966, 1386, 8352),
    @(44253, 7075, 1420, 8495),
    @(44254, 7189, 1441, 8630),
    @(44255, 7270, 1472, 8742),
    @(44256, 7388, 1497, 8885),
    @(44257, 7489, 1524, 9013),
    @(44258, 7560, 1545, 9105),
    @(44259, 7665, 1559, 9224),
    @(44260, 7739, 1578, 9317),
    @(44261, 7836, 1589, 9425),
    @(44262, 7921, 1594, 9515)
)

$startRow = 123
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
}

$lastRow = $startRow + $newRows.Count - 1

# Scroll the view down and select the newly added last cell (A143), matching
# the saved workbook's sheet view state after the edit.
$excel.ActiveWindow.ScrollRow = 97
$ws.Range("A$lastRow").Select()
